# Auto-applied edits matching the authoritative diff (crypto price/volume refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format cells whose new price text is numeric-looking as Text,
# so Excel stores the literal digits/dots instead of coercing to a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the refreshed Price (D) and Volume(1h) (E) text for each row.
$ws.Range("D2").Value = "68.655.25"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "2.454.53"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "557.45"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "161.40"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("D9").Value = "0.152"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  +0.66%  "
$ws.Range("D11").Value = "0.331"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("D13").Value = "68.572.18"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -1.60%  "
$ws.Range("D15").Value = "23.34"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "10.56"
$ws.Range("E16").Value = "  -3.73%  "
$ws.Range("D17").Value = "334.48"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").Value = "6.90"
$ws.Range("E18").Value = "  -3.73%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").Value = "66.25"
$ws.Range("E22").Value = "  -2.17%  "
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "8.15"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").Value = "0.0₃0813"
$ws.Range("E25").Value = "  -2.54%  "
$ws.Range("D26").Value = "7.15"
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "425.40"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  -3.23%  "
$ws.Range("D30").Value = "1.61"
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("D31").Value = "158.46"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("D32").Value = "19.01"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  -1.62%  "
$ws.Range("D35").Value = "17.72"
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("D37").Value = "4.38"
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  -4.09%  "
$ws.Range("E39").Value = "  -2.32%  "
$ws.Range("D40").Value = "2.04"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("E41").Value = "  -0.65%  "
$ws.Range("D42").Value = "129.12"
$ws.Range("E42").Value = "  -3.39%  "
$ws.Range("D43").Value = "0.0717"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "0.479"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.0911"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  -3.60%  "
$ws.Range("D49").Value = "4.91"
$ws.Range("E49").Value = "  -8.28%  "
$ws.Range("D50").Value = "16.69"
$ws.Range("E50").Value = "  -4.85%  "
$ws.Range("E51").Value = "  -3.21%  "
